$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 5).Value = "NA"
}

# Restore F1 (unrelated cell) to blank; some runtimes materialize an
# empty shared-string placeholder cell into a real value on save.
$ws.Range("F1").Value = $null
